$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Internal Assignment" column (column O) contents: header + data rows
$ws.Range("O4:O8").ClearContents()

# Update the active cell selection to O6 (matches the recorded selection after the edit)
$ws.Range("O6").Select()
